# Commit: "Try to fix statement on pt 9"
# The paragraph discussing compliance with "п. 9 Положения о присуждении ученых
# степеней" is rewritten: a new explanatory sentence is inserted, and the
# description of the algorithms is trimmed down while gaining a couple of
# extra qualifying words ("точный", "и построения нижней оценки и").

$d = $word.ActiveDocument

# 1) Insert the new explanatory text right after "...исследований" and before
#    the (still yellow-highlighted) "разработаны алгоритм..." text, also
#    inserting "точный " right after the (now duplicated) "разработаны " and
#    before "алгоритм".
$d.Content.Find.Execute(
    "исследований разработаны алгоритм",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "исследований разработаны теоретические положения, совокупность которых можно квалифицировать как научное достижение. В диссертации содержится решение научной задачи, имеющей важное значение для развития соответствующей отрасли знаний (???), а именно разработаны точный алгоритм",
    2
) | Out-Null

# 2) Trim the long tail describing the algorithms down to the shorter phrase,
#    adding "и построения нижней оценки и" before the heuristic algorithm
#    mention, and dropping the rest of the old sentence.
$d.Content.Find.Execute(
    "с ограничениями предшествования, эвристический алгоритм решения задачи непрерывной резки, схемы информационного обмена и методика использования разработанных алгоритмов в системах автоматизированного проектирования управляющих программ машин листовой резки с ЧПУ, имеющие существенное значение при оптимизации технологических процессов раскройно-заготовительного производства в машиностроении и других отраслях промышленности Российской Федерации.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "с ограничениями предшествования и построения нижней оценки и эвристический алгоритм решения задачи непрерывной резки.",
    2
) | Out-Null

# 3) The second "разработаны ... резки" phrase (and everything up through it)
#    keeps the original yellow highlight; re-apply it explicitly so the newly
#    inserted words ("точный", "и построения нижней оценки и") are covered too.
$hiRange = $d.Content
$hiRange.Find.Execute(
    "разработаны точный алгоритм ветвей и границ для решения обобщенной задачи коммивояжера с ограничениями предшествования и построения нижней оценки и эвристический алгоритм решения задачи непрерывной резки",
    $true, $false, $false, $false, $false, $true, 1, $false
) | Out-Null
$hiRange.HighlightColorIndex = 7
